$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting the existing
# "dbExcel" / "WebExcel" columns one slot to the right.
$ws.Columns("B:B").Insert()

# Match the width of the newly inserted column to column A.
$ws.Columns("B:B").ColumnWidth = 75

# New header + query text for the inserted "StatQuery" column.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Squamous cell lung carcinoma']  OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match the wrap-text style used by the existing query cell (A2).
$ws.Range("B2").WrapText = $true

# Update the selection to reflect the saved view state.
$ws.Range("A2").Select() | Out-Null
